$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the old, never-used empty placeholder rows (4 through 26).
#    Only rows 2 and 3 held real data; everything below was just a
#    leftover template row with a stray format on column G.
# ------------------------------------------------------------------
$ws.Range("4:26").Delete()

# ------------------------------------------------------------------
# 2. Re-create row 4 (subject 3) and row 5 (subject 4) by cloning the
#    formatting of the existing data rows, then filling in the values.
# ------------------------------------------------------------------
$ws.Range("B3:H3").Copy()
$ws.Range("B4:H4").PasteSpecial(-4122)

$ws.Range("B2:E2").Copy()
$ws.Range("B5:E5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Subject 3 - fully completed entry
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 26
$ws.Range("D4").Value = "F"
$ws.Range("E4").Value = 0.61319444444444449
$ws.Range("F4").Value = 0.62638888888888888
$ws.Range("G4").Value = 12
$ws.Range("H4").Value = "when starting tutorial (no feedback), she was expecting to see the red version first (it is the first time the black one is shown first)  + problems visualizing the text with Safari (version 14.1.1 ) + it's written ""right or left"" arrow, but instead it should be ""left or right"" + she did not feel tired at all, could have done other 2/3 blocks"
$ws.Rows(4).RowHeight = 75

# Subject 4 - link sent only so far
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 27
$ws.Range("D5").Value = "M"
$ws.Range("E5").Value = 0.60833333333333328
$ws.Rows(5).RowHeight = 24

# ------------------------------------------------------------------
# 3. Update the comments column: subject 1 learned a new detail
#    (Safari version) and the comment text box alignment for the
#    whole column is unified to center/middle + wrap.
# ------------------------------------------------------------------
$ws.Range("H2").Value = $ws.Range("H2").Value2 + " + safari version 16.4"

$ws.Range("H2:H4").HorizontalAlignment = -4108   # xlCenter
$ws.Range("H2:H4").VerticalAlignment = -4108     # xlCenter
$ws.Range("H2:H4").WrapText = $true
$ws.Range("H2:H4").IndentLevel = 0

$ws.Rows(2).RowHeight = 49.8

# ------------------------------------------------------------------
# 4. Tutorial duration column now shows two decimals without the
#    thousands separator.
# ------------------------------------------------------------------
$ws.Range("G1:G4").NumberFormat = "0.00"

# ------------------------------------------------------------------
# 5. Match the final selection left behind in the source file.
# ------------------------------------------------------------------
$ws.Range("E5").Select()
